# "sn from new process" - append newly-processed key/lookup rows to the
# "key" worksheet (ego/alter organization match workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("key")
$ws.Activate()

# New rows to append below the existing key table (rows 2-216).
# Columns: A = raw name, B = normalized/matched name, C/D = category.
$newRows = @(
    @("Commercial Diver 1", "Commercial Diver 1", "Individual", "Individual"),
    @("Commercial Diver 4", "Commercial Diver 4", "Individual", "Individual"),
    @("Commercial Diver 8", "Commercial Diver 8", "Individual", "Individual"),
    @("Commercial Diver 2", "Commercial Diver 2", "Individual", "Individual"),
    @("Commercial Diver 5", "Commercial Diver 5", "Individual", "Individual"),
    @("Commercial Diver 6", "Commercial Diver 6", "Individual", "Individual"),
    @("Commercial Diver 7", "Commercial Diver 7", "Individual", "Individual"),
    @("Commercial Diver 3", "Commercial Diver 3", "Individual", "Individual"),
    @("Photographer 1",     "Artist 2",           "Individual", "Individual"),
    @("Artist 1",           "Artist 1",           "Individual", "Individual")
)

$startRow = 217
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

# Mirror the author's final scroll/selection state in the saved view.
$excel.ActiveWindow.ScrollRow = 199
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A230").Select() | Out-Null
